$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.270.76"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.688.75"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "219.18"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "0.5248"
$ws.Range("E6").Value = "  +4.62%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").Value = "0.06431"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").Value = "22.08"
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("D11").Value = "0.07467"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "1.696.79"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "4.552"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "0.5862"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "0.000008551"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "64.65"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "26.315.30"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "4.975"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "190.46"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").Value = "6.242"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "145.27"
$ws.Range("D25").Value = "7.682"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  +6.57%  "
$ws.Range("D27").Value = "15.85"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").Value = "0.06717"
$ws.Range("E28").Value = "  +15.82%  "
$ws.Range("D29").Value = "1.347"
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").Value = "1.331"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").Value = "3.602"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("D32").Value = "3.556"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "1.667"
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").Value = "1.029"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").Value = "0.6211"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("D36").Value = "2.386"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("D38").Value = "6.287"
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").Value = "0.01619"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").Value = "1.103.22"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "0.8774"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("D42").Value = "1.017"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "100.92"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").Value = "1.838.83"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "0.00000000116"
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("D46").Value = "56.91"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").Value = "8.151"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "0.05263"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "0.4299"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "6.019"
$ws.Range("E51").Value = "  +3.33%  "
